# Apply updated cryptocurrency price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.383.00"
$ws.Range("E2").Value = "  -4.36%  "

$ws.Range("D3").Value = "1.568.45"
$ws.Range("E3").Value = "  -4.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3680"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3393"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.175"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07598"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.047"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.901"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001138"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.94%  "

$ws.Range("D17").Value = "1.569.51"
$ws.Range("E17").Value = "  -4.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06770"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.91%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.231"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5314"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.56%  "

$ws.Range("D25").Value = "22.390.83"
$ws.Range("E25").Value = "  -4.40%  "

$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.983"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.43%  "

$ws.Range("E29").Value = "  -4.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.978"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.18%  "

$ws.Range("D32").Value = "1.742.06"
$ws.Range("E32").Value = "  -4.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.044"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.254"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.44%  "

$ws.Range("E35").Value = "  -6.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08445"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02545"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2326"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06554"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.531"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.249"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6373"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.11%  "

$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.778"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.135"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.47%  "
